# Update cryptocurrency Price (col D) and Volume(1h) (col E) values
# for rows 2-51 of the active sheet, matching the refreshed data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values look like plain decimal numbers (e.g. "1.005").
# Force those cells to text format first so Excel keeps them as strings
# instead of silently converting them to numeric values.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.459.79"
$ws.Range("E2").Value = "  -3.65%  "
$ws.Range("D3").Value = "1.861.11"
$ws.Range("E3").Value = "  -4.62%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("D5").Value = "321.44"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").Value = "0.4502"
$ws.Range("E7").Value = "  -5.71%  "
$ws.Range("D8").Value = "0.3861"
$ws.Range("E8").Value = "  -4.04%  "
$ws.Range("D9").Value = "47.35"
$ws.Range("E9").Value = "  -11.66%  "
$ws.Range("D10").Value = "0.08010"
$ws.Range("E10").Value = "  -5.95%  "
$ws.Range("D11").Value = "1.020"
$ws.Range("E11").Value = "  -3.94%  "
$ws.Range("D12").Value = "21.46"
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("D13").Value = "1.878.21"
$ws.Range("E13").Value = "  -4.06%  "
$ws.Range("D14").Value = "5.887"
$ws.Range("E14").Value = "  -4.88%  "
$ws.Range("D15").Value = "7.139"
$ws.Range("E15").Value = "  -6.34%  "
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").Value = "0.06547"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").Value = "17.20"
$ws.Range("E20").Value = "  -8.18%  "
$ws.Range("D21").Value = "1.006"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").Value = "5.526"
$ws.Range("E22").Value = "  -4.93%  "
$ws.Range("D23").Value = "27.472.16"
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("D24").Value = "10.86"
$ws.Range("E24").Value = "  -5.75%  "
$ws.Range("D25").Value = "2.281"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").Value = "2.099.81"
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("D27").Value = "151.23"
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("D28").Value = "19.46"
$ws.Range("E28").Value = "  -3.54%  "
$ws.Range("D29").Value = "5.534"
$ws.Range("E29").Value = "  -7.02%  "
$ws.Range("D30").Value = "2.038"
$ws.Range("E30").Value = "  -5.71%  "
$ws.Range("D31").Value = "121.06"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "0.09392"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").Value = "1.485"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("D34").Value = "0.9321"
$ws.Range("E34").Value = "  -6.37%  "
$ws.Range("D35").Value = "3.622"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "5.293"
$ws.Range("E36").Value = "  -5.51%  "
$ws.Range("D37").Value = "0.02234"
$ws.Range("E37").Value = "  -4.61%  "
$ws.Range("D38").Value = "1.228"
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("D39").Value = "0.05972"
$ws.Range("E39").Value = "  -4.17%  "
$ws.Range("D40").Value = "8.359"
$ws.Range("E40").Value = "  -4.78%  "
$ws.Range("D41").Value = "1.005"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("E42").Value = "  -4.62%  "
$ws.Range("D43").Value = "0.1859"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("E44").Value = "  -7.30%  "
$ws.Range("D45").Value = "1.277"
$ws.Range("E45").Value = "  -3.78%  "
$ws.Range("D46").Value = "0.5674"
$ws.Range("E46").Value = "  -4.93%  "
$ws.Range("D47").Value = "12.16"
$ws.Range("E47").Value = "  -6.23%  "
$ws.Range("D48").Value = "1.934"
$ws.Range("E48").Value = "  -6.31%  "
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "1.005"
$ws.Range("E51").Value = "  -11.15%  "
